# Case 5_156 (380 kV): slack bus voltage setpoint changed from 1.05 pu to 1.02 pu.
# Re-run load-flow results (res_bus/vm_pu) for rows 2-25, columns B-F and I-M.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025109979222057
$ws.Range("D2").Value = 1.028327507728955
$ws.Range("E2").Value = 1.025421719044825
$ws.Range("F2").Value = 1.023603807099087
$ws.Range("I2").Value = 1.029356175659613
$ws.Range("J2").Value = 1.030280834624214
$ws.Range("K2").Value = 1.031144646376276
$ws.Range("L2").Value = 1.028247338075449
$ws.Range("M2").Value = 1.026434758240348
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026816996821466
$ws.Range("D3").Value = 1.029572503769889
$ws.Range("E3").Value = 1.026897240026056
$ws.Range("F3").Value = 1.025958908886571
$ws.Range("I3").Value = 1.029763678090753
$ws.Range("J3").Value = 1.031623137731457
$ws.Range("K3").Value = 1.032196375637372
$ws.Range("L3").Value = 1.029528345412348
$ws.Range("M3").Value = 1.02859256105393
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027917597803258
$ws.Range("D4").Value = 1.030374573288738
$ws.Range("E4").Value = 1.027848682398705
$ws.Range("F4").Value = 1.027478105170414
$ws.Range("I4").Value = 1.030024036439695
$ws.Range("J4").Value = 1.0324875167029
$ws.Range("K4").Value = 1.032872847853746
$ws.Range("L4").Value = 1.030353437239987
$ws.Range("M4").Value = 1.029983813633798
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028379362758668
$ws.Range("D5").Value = 1.030710931655675
$ws.Range("E5").Value = 1.028247889577853
$ws.Range("F5").Value = 1.028115678875713
$ws.Range("I5").Value = 1.030132700974056
$ws.Range("J5").Value = 1.032849916000458
$ws.Range("K5").Value = 1.033156274971455
$ws.Range("L5").Value = 1.030699409022867
$ws.Range("M5").Value = 1.030567530743169
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028456841267622
$ws.Range("D6").Value = 1.030767359297107
$ws.Range("E6").Value = 1.028314872966502
$ws.Range("F6").Value = 1.02822266692482
$ws.Range("I6").Value = 1.030150900050432
$ws.Range("J6").Value = 1.032910707080482
$ws.Range("K6").Value = 1.033203807571908
$ws.Range("L6").Value = 1.030757447027979
$ws.Range("M6").Value = 1.030665471911107
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027923771548911
$ws.Range("D7").Value = 1.030379070976491
$ws.Range("E7").Value = 1.027854019665564
$ws.Range("F7").Value = 1.027486628726976
$ws.Range("I7").Value = 1.030025491517732
$ws.Range("J7").Value = 1.032492362952099
$ws.Range("K7").Value = 1.032876638783877
$ws.Range("L7").Value = 1.030358063636423
$ws.Range("M7").Value = 1.029991617830588
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025687703796404
$ws.Range("D8").Value = 1.028748997667068
$ws.Range("E8").Value = 1.025921074579339
$ws.Range("F8").Value = 1.024400719605939
$ws.Range("I8").Value = 1.029494584338202
$ws.Range("J8").Value = 1.030735347187083
$ws.Range("K8").Value = 1.031500932459475
$ws.Range("L8").Value = 1.028681057095857
$ws.Range("M8").Value = 1.027165048704868
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.021716337316919
$ws.Range("D9").Value = 1.025849027034356
$ws.Range("E9").Value = 1.022488874597086
$ws.Range("F9").Value = 1.018925432794981
$ws.Range("I9").Value = 1.028533362337054
$ws.Range("J9").Value = 1.027606548005155
$ws.Range("K9").Value = 1.029045078098312
$ws.Range("L9").Value = 1.025696166649526
$ws.Range("M9").Value = 1.02214473490405
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019046632427948
$ws.Range("D10").Value = 1.023896373917577
$ws.Range("E10").Value = 1.020182219179631
$ws.Range("F10").Value = 1.015247998947158
$ws.Range("I10").Value = 1.027874923260944
$ws.Range("J10").Value = 1.025497711254171
$ws.Range("K10").Value = 1.027385782129007
$ws.Range("L10").Value = 1.023685301407006
$ws.Range("M10").Value = 1.01876941823197
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.0178851033332
$ws.Range("D11").Value = 1.023046086031799
$ws.Range("E11").Value = 1.019178800278517
$ws.Range("F11").Value = 1.013648715430646
$ws.Range("I11").Value = 1.027585552823296
$ws.Range("J11").Value = 1.024578893659922
$ws.Range("K11").Value = 1.026661885708661
$ws.Range("L11").Value = 1.022809403873174
$ws.Range("M11").Value = 1.017300705986151
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017452805628249
$ws.Range("D12").Value = 1.022729517949614
$ws.Range("E12").Value = 1.018805372446812
$ws.Range("F12").Value = 1.013053589902746
$ws.Range("I12").Value = 1.027477420534274
$ws.Range("J12").Value = 1.024236731824292
$ws.Range("K12").Value = 1.026392170625026
$ws.Range("L12").Value = 1.022483260572662
$ws.Range("M12").Value = 1.016754046313995
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017545573981414
$ws.Range("D13").Value = 1.02279745631656
$ws.Range("E13").Value = 1.018885506573952
$ws.Range("F13").Value = 1.013181295876004
$ws.Range("I13").Value = 1.027500644677788
$ws.Range("J13").Value = 1.024310166472512
$ws.Range("K13").Value = 1.026450063123629
$ws.Range("L13").Value = 1.022553255724774
$ws.Range("M13").Value = 1.016871357691545
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.017849387029558
$ws.Range("D14").Value = 1.023019933466438
$ws.Range("E14").Value = 1.019147947280976
$ws.Range("F14").Value = 1.013599544438194
$ws.Range("I14").Value = 1.027576627818142
$ws.Range("J14").Value = 1.024550628357462
$ws.Range("K14").Value = 1.026639607947942
$ws.Range("L14").Value = 1.022782461107331
$ws.Range("M14").Value = 1.017255541874831
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018036462396112
$ws.Range("D15").Value = 1.023156911336475
$ws.Range("E15").Value = 1.01930955056009
$ws.Range("F15").Value = 1.013857096831287
$ws.Range("I15").Value = 1.027623357566377
$ws.Range("J15").Value = 1.024698668726322
$ws.Range("K15").Value = 1.026756282639454
$ws.Range("L15").Value = 1.022923576099638
$ws.Range("M15").Value = 1.017492101649632
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.019123600806407
$ws.Range("D16").Value = 1.023952702758047
$ws.Range("E16").Value = 1.020248713795414
$ws.Range("F16").Value = 1.015353988778982
$ws.Range("I16").Value = 1.027894037441805
$ws.Range("J16").Value = 1.025558568911837
$ws.Range("K16").Value = 1.027433709496295
$ws.Range("L16").Value = 1.023743321201805
$ws.Range("M16").Value = 1.018866737524127
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019804037279295
$ws.Range("D17").Value = 1.024450591940156
$ws.Range("E17").Value = 1.02083657535601
$ws.Range("F17").Value = 1.016291066600944
$ws.Range("I17").Value = 1.028062681883372
$ws.Range("J17").Value = 1.02609642818957
$ws.Range("K17").Value = 1.027857182538998
$ws.Range("L17").Value = 1.024256126451415
$ws.Range("M17").Value = 1.019727063227277
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020200392464898
$ws.Range("D18").Value = 1.024740542518065
$ws.Range("E18").Value = 1.021179020410994
$ws.Range("F18").Value = 1.016836982236367
$ws.Range("I18").Value = 1.028160638587447
$ws.Range("J18").Value = 1.026409605928965
$ws.Range("K18").Value = 1.02810366620078
$ws.Range("L18").Value = 1.024554738534573
$ws.Range("M18").Value = 1.020228186565726
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020335449823456
$ws.Range("D19").Value = 1.024839330661269
$ws.Range("E19").Value = 1.021295710475441
$ws.Range("F19").Value = 1.017023013571992
$ws.Range("I19").Value = 1.02819396982979
$ws.Range("J19").Value = 1.026516299326577
$ws.Range("K19").Value = 1.028187622835848
$ws.Range("L19").Value = 1.024656473510397
$ws.Range("M19").Value = 1.020398940708821
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019731088070227
$ws.Range("D20").Value = 1.024397220811134
$ws.Range("E20").Value = 1.020773549502739
$ws.Range("F20").Value = 1.016190596189034
$ws.Range("I20").Value = 1.028044630465902
$ws.Range("J20").Value = 1.026038777641628
$ws.Range("K20").Value = 1.027811801885784
$ws.Range("L20").Value = 1.024201158995921
$ws.Range("M20").Value = 1.0196348300606
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017759945425523
$ws.Range("D21").Value = 1.022954439861219
$ws.Range("E21").Value = 1.019070684880479
$ws.Range("F21").Value = 1.013476410804059
$ws.Range("I21").Value = 1.027554270599782
$ws.Range("J21").Value = 1.024479842554068
$ws.Range("K21").Value = 1.026583814691333
$ws.Range("L21").Value = 1.0227149879514
$ws.Range("M21").Value = 1.017142440268609
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016515656519328
$ws.Range("D22").Value = 1.022043055576406
$ws.Range("E22").Value = 1.01799588932371
$ws.Range("F22").Value = 1.011763623867019
$ws.Range("I22").Value = 1.027242213739954
$ws.Range("J22").Value = 1.023494622974745
$ws.Range("K22").Value = 1.02580693391507
$ws.Range("L22").Value = 1.021775959062057
$ws.Range("M22").Value = 1.015568909117158
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017175755021655
$ws.Range("D23").Value = 1.022526605859108
$ws.Range("E23").Value = 1.018566057337284
$ws.Range("F23").Value = 1.012672212286494
$ws.Range("I23").Value = 1.02740799860509
$ws.Range("J23").Value = 1.024017392148712
$ws.Range("K23").Value = 1.026219232796392
$ws.Range("L23").Value = 1.022274199419794
$ws.Range("M23").Value = 1.016403692948373
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019764052320646
$ws.Range("D24").Value = 1.024421338347303
$ws.Range("E24").Value = 1.020802029546396
$ws.Range("F24").Value = 1.016235996499964
$ws.Range("I24").Value = 1.028052788393474
$ws.Range("J24").Value = 1.026064829130785
$ws.Range("K24").Value = 1.027832309059089
$ws.Range("L24").Value = 1.024225997963028
$ws.Range("M24").Value = 1.019676508389528
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.022746841975043
$ws.Range("D25").Value = 1.026602088188735
$ws.Range("E25").Value = 1.023379373133355
$ws.Range("F25").Value = 1.020345577786254
$ws.Range("I25").Value = 1.028784941805884
$ws.Range("J25").Value = 1.02841939385682
$ws.Range("K25").Value = 1.029683807264326
$ws.Range("L25").Value = 1.026471455458504
$ws.Range("M25").Value = 1.023447481378648
